$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three renamed font samples and one source link
# (order matters: new shared-string entries are appended in the order
# the values are assigned)
$ws.Cells.Item(19, 2).Value = "Dromedar"
$ws.Cells.Item(23, 2).Value = "Zeppelin"
$ws.Cells.Item(12, 2).Value = "Pacman"
$ws.Cells.Item(23, 5).Value = "Wikipedia"

# Update the view state (scrolled/selected cell) to match the new window
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("E30").Select()
